# aggiornamento fino a 20/09/2021
# Append new daily rows (375-385) to the single worksheet, extending the
# existing data table with the same layout as the prior rows:
#   A = date serial (styled like the rest of column A)
#   B = nuovi pos. (new positives)
#   C = somma mobile 7gg. (7-day rolling sum)
#   D = computed rolling-average-derived value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 375
$lastNewRow  = 385

# Carry the date-column formatting (style used by A2:A374) down into the
# newly appended date cells, without introducing new style/font entries.
$ws.Range("A374").Copy() | Out-Null
$destRange = $ws.Range($ws.Cells.Item($firstNewRow, 1), $ws.Cells.Item($lastNewRow, 1))
$destRange.PasteSpecial(-4122) | Out-Null

$data = @(
    @(44449, 0, 4, 56.89900426742533),
    @(44450, 0, 4, 56.89900426742533),
    @(44451, 3, 7, 99.5732574679943),
    @(44452, 1, 7, 99.5732574679943),
    @(44453, 0, 5, 71.12375533428165),
    @(44454, 0, 5, 71.12375533428165),
    @(44455, 0, 4, 56.89900426742533),
    @(44456, 2, 6, 85.34850640113798),
    @(44457, 0, 6, 85.34850640113798),
    @(44458, 1, 4, 56.89900426742533),
    @(44459, 0, 3, 42.67425320056899)
)

$r = $firstNewRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
